$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7234810590744019
$ws.Range("B1").Value = 2.11238956451416
$ws.Range("C1").Value = 6.327485084533691
$ws.Range("D1").Value = 1.77379834651947
$ws.Range("E1").Value = 1.050476431846619
